$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column D contain price strings such as "63.988.19" that must
# remain plain text (they are not valid numbers). Force the column to
# a text format before writing so Excel does not coerce them into
# floating point numbers (which would lose precision/trailing zeros),
# then restore the original "Normal" style so cell styling is unaffected.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.944.62"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.637.89"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "579.37"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").Value = "156.79"
$ws.Range("E6").Value = "  +0.62%  "

$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  -2.25%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "2.635.99"
$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("E10").Value = "  -2.59%  "

$ws.Range("E11").Value = "  +0.42%  "

$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").Value = "28.67"
$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("D15").Value = "3.114.74"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "63.867.90"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "2.636.36"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "7.76"
$ws.Range("E20").Value = "  +2.61%  "

$ws.Range("E21").Value = "  -2.74%  "

$ws.Range("D22").Value = "345.01"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "68.30"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("E25").Value = "  +7.75%  "

$ws.Range("E26").Value = "  +2.89%  "

$ws.Range("D27").Value = "9.28"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("E28").Value = "  +3.66%  "

$ws.Range("D29").Value = "581.19"
$ws.Range("E29").Value = "  +1.29%  "

$ws.Range("E30").Value = "  +3.34%  "

$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("E35").Value = "  +2.45%  "

$ws.Range("E36").Value = "  +2.78%  "

$ws.Range("E37").Value = "  -1.65%  "

$ws.Range("D38").Value = "19.77"
$ws.Range("E38").Value = "  -0.86%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("D41").Value = "153.29"
$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("E42").Value = "  +7.78%  "

$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").Value = "162.24"
$ws.Range("E44").Value = "  +3.25%  "

$ws.Range("D45").Value = "24.18"
$ws.Range("E45").Value = "  +4.31%  "

$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").Value = "0.0590"
$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("D48").Value = "0.636"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("E49").Value = "  -2.05%  "

$ws.Range("E50").Value = "  -1.34%  "

$ws.Range("D51").Value = "19.07"
$ws.Range("E51").Value = "  -0.16%  "

# Restore default styling on column D (clears the temporary text format).
$priceRange.Style = "Normal"
